$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the D (Price) column keeps values as plain text, matching the
# original inline-string cell type, even for values that look numeric.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.920.47'
$ws.Range("E2").Value = '  -4.35%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.608.38'
$ws.Range("E3").Value = '  -4.23%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.49'
$ws.Range("E5").Value = '  -1.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.80'
$ws.Range("E6").Value = '  -2.26%  '

$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").Value = '  -1.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.70'
$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("E10").Value = '  -3.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.335'
$ws.Range("E11").Value = '  -1.13%  '

$ws.Range("E12").Value = '  +1.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.075.91'
$ws.Range("E13").Value = '  -3.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '57.974.94'
$ws.Range("E14").Value = '  -4.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.58'
$ws.Range("E15").Value = '  -3.16%  '

$ws.Range("E16").Value = '  -2.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.624.04'
$ws.Range("E17").Value = '  -3.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.39'
$ws.Range("E18").Value = '  -2.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '333.97'
$ws.Range("E19").Value = '  -3.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.29'
$ws.Range("E20").Value = '  -3.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.24'
$ws.Range("E21").Value = '  -3.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.80'
$ws.Range("E23").Value = '  +0.44%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.415'
$ws.Range("E24").Value = '  -1.44%  '

$ws.Range("E25").Value = '  -2.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.03'
$ws.Range("E27").Value = '  -3.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0780'
$ws.Range("E28").Value = '  -4.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.58'
$ws.Range("E29").Value = '  -3.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  -1.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '150.76'
$ws.Range("E32").Value = '  +0.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.64'
$ws.Range("E33").Value = '  -2.12%  '

$ws.Range("E34").Value = '  -4.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.16'
$ws.Range("E35").Value = '  -5.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.887'
$ws.Range("E36").Value = '  -5.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.49'
$ws.Range("E37").Value = '  -1.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.838'
$ws.Range("E38").Value = '  -3.80%  '

$ws.Range("E39").Value = '  -6.17%  '

$ws.Range("E40").Value = '  -2.00%  '

$ws.Range("E41").Value = '  +0.18%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0963'
$ws.Range("E42").Value = '  -2.26%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.594'
$ws.Range("E43").Value = '  -2.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '266.66'
$ws.Range("E44").Value = '  -5.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.60'
$ws.Range("E45").Value = '  +1.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.04'
$ws.Range("E46").Value = '  -5.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0529'
$ws.Range("E47").Value = '  -1.68%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.028.20'
$ws.Range("E48").Value = '  -5.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0226'
$ws.Range("E49").Value = '  -2.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.60'
$ws.Range("E50").Value = '  -4.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.14'
$ws.Range("E51").Value = '  -4.96%  '
